# Adds two new columns, I ("I0") and J ("IF"), to Sheet1, mirroring the
# header style already used by the other header cells (B1:H1) and filling
# in the numeric data for rows 2-73.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells -----------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Clone the header formatting (bold font, thin border, centered/top
# aligned) from the existing "IP" header cell (H1) onto the two new
# header cells so they visually match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Body data (rows 2-73) --------------------------------------------
$i0Values = @(9, 6, 7, 6, 7, 5, 7, 5, 10, 7, 8, 7, 6, 8, 6, 7, 6, 7, 8, 8, 6, 7, 8, 7, 7, 7, 8, 9, 9, 6, 7, 8, 8, 7, 8, 8, 8, 8, 7, 7, 8, 6, 7, 7, 8, 9, 7, 7, 7, 6, 7, 6, 7, 8, 8, 8, 7, 8, 9, 8, 8, 8, 7, 4, 1, 5, 8, 6, 5, 5, 5, 5)
$ifValues = @(9, 7, 7, 6, 7, 5, 7, 6, 10, 7, 8, 7, 7, 8, 6, 7, 6, 7, 8, 8, 6, 7, 8, 7, 7, 7, 8, 9, 9, 6, 7, 8, 8, 7, 8, 8, 8, 8, 7, 7, 8, 6, 7, 7, 8, 9, 7, 7, 7, 6, 7, 6, 7, 8, 8, 8, 7, 8, 9, 8, 8, 8, 7, 4, 1, 5, 8, 6, 5, 5, 5, 5)

for ($idx = 0; $idx -lt $i0Values.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $i0Values[$idx]
    $ws.Cells.Item($row, 10).Value = $ifValues[$idx]
}

Write-Host "Added I0/IF columns (I1:J73)"
